$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.402.52"
$ws.Range("D3").Value = "1.843.42"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'0.6319"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.07529"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "'0.2926"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "'24.42"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'0.07715"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.842.89"
$ws.Range("E12").Value = "  -7.16%  "
$ws.Range("D14").Value = "'0.6792"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D16").Value = "'83.18"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "2.096.63"
$ws.Range("E17").Value = "  -7.41%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "29.429.99"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'228.27"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'7.450"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'157.18"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "'0.1393"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'8.375"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'1.459"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'1.281"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").Value = "'0.05631"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").Value = "'4.100"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "'4.030"
$ws.Range("D34").Value = "'1.840"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.7097"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "1.245.68"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'2.766"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'6.333"
$ws.Range("E41").Value = "  +3.95%  "
$ws.Range("D42").Value = "'0.9017"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'101.74"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'65.80"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "'0.00000000118"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'7.098"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'0.3997"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.920"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.672"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -0.46%  "
